# Rename all 30 worksheets from their old "summ<number>" names to the new
# "summ<number>" names, keeping their order, sheetId, and relationship ids
# unchanged (only the visible sheet name changes), per the commit
# "rerun models and create results figures and tables".

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ54439426",
    "summ54592918",
    "summ54778215",
    "summ54965360",
    "summ55137465",
    "summ55316133",
    "summ55481244",
    "summ55736464",
    "summ55875546",
    "summ56016960",
    "summ56152238",
    "summ56296183",
    "summ56442315",
    "summ56584702",
    "summ56735380",
    "summ56881422",
    "summ57020409",
    "summ57159972",
    "summ57303039",
    "summ57435734",
    "summ57580010",
    "summ57734301",
    "summ57877479",
    "summ58005711",
    "summ58150354",
    "summ58293304",
    "summ58475967",
    "summ58618448",
    "summ58755876",
    "summ58899266"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
